$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: original rows (language, value) excluding zero-value entries
# (Russian, Uzbek), sorted descending by value.
$data = @(
    @("English", 30.06594013551159),
    @("Spanish", 8.92153640356304),
    @("Japanese", 8.606511260446988),
    @("German", 7.295965501115992),
    @("Arabic", 5.202772886067171),
    @("Portuguese", 4.630497098882373),
    @("Chinese", 4.434789486162071),
    @("French", 4.413647053662403),
    @("Italian", 4.314254429945956),
    @("Malay-Indonesian", 2.041355970101624),
    @("Dutch", 1.877656724682577),
    @("Persian", 1.583044464864306),
    @("Turkish", 1.364097569620373),
    @("Polish", 1.054763876569358),
    @("Korean", 0.9200050266423154),
    @("Urdu", 0.7369751430991435),
    @("Swedish", 0.6629717978986374),
    @("Thai", 0.6536760848891855),
    @("Bengali", 0.4110164248581473),
    @("Vietnamese", 0.2690055659772652)
)

$lastOldRow = 23
$newLastRow = 1 + $data.Count

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove any leftover rows beyond the new data range (previously rows 22-23)
if ($newLastRow -lt $lastOldRow) {
    $startClear = $newLastRow + 1
    $clearRange = $ws.Range($ws.Cells.Item($startClear, 1), $ws.Cells.Item($lastOldRow, 2))
    $clearRange.Delete()
}
